$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update password cell value
$ws.Range("B2").Value = "kanbas123$"

# Update active cell selection
$ws.Range("B2").Select()
